# Invoice 4.xlsx edit script
# - Update the "CLIENT CODE" placeholder text to a concrete name/email.
# - Update the "Automatic Door" quantity/amount from 1 to 120 (recalculates
#   SUBTOTAL, TAX and TOTAL automatically via existing formulas).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Invoice")

# Update client code placeholder text
$ws.Range("A31").Value = "Charlie, charlie@mail.com"

# Update the Automatic Door amount (E18) -- dependent formulas (F18, F21,
# F23, F24) recalc automatically.
$ws.Range("E18").Value = 120
